# Apply "Update gh-pages to output generated at 456a3b4" changes.
# This updates the "想去人数" (interested-count) column F on three sheets
# (展览, 演出, 全部类型) and appends one new row to 本地生活 (with a matching
# update of its dimension), mirroring the source XML diff.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "展览" (Exhibitions) - column F increments
# ---------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    3  = 969
    7  = 1190
    8  = 948
    9  = 36
    11 = 1051
    12 = 2513
    13 = 576
    15 = 1697
    17 = 648
    18 = 20
    22 = 1528
    23 = 772
    24 = 657
    25 = 515
    28 = 45
    30 = 1165
    31 = 333
    32 = 2462
    35 = 470
    38 = 4036
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Cells.Item($row, 6).Value = $expoUpdates[$row]
}

# ---------------------------------------------------------------------
# Sheet "演出" (Performances) - column F increments
# ---------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$showUpdates = @{
    7  = 644
    14 = 4142
    23 = 263
    25 = 129
    28 = 45
}
foreach ($row in $showUpdates.Keys) {
    $wsShow.Cells.Item($row, 6).Value = $showUpdates[$row]
}

# ---------------------------------------------------------------------
# Sheet "本地生活" (Local life) - column F increments + new row 8
# ---------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$localUpdates = @{
    5 = 1683
    6 = 459
    7 = 1035
}
foreach ($row in $localUpdates.Keys) {
    $wsLocal.Cells.Item($row, 6).Value = $localUpdates[$row]
}

# Append the new event row (row 8)
$wsLocal.Cells.Item(8, 1).Value = 7
$wsLocal.Cells.Item(8, 2).Value = "'2024-03-21"
$wsLocal.Cells.Item(8, 3).Value = "上海·NIJISANJI EN 官方授权主题店"
$wsLocal.Cells.Item(8, 4).Value = "西藏北路166号（地铁8号线曲阜路下） 静安大悦城"
$wsLocal.Cells.Item(8, 5).Value = "2024.03.21 00:00-04.28 23:59"
$wsLocal.Cells.Item(8, 6).Value = 33
$wsLocal.Cells.Item(8, 7).Value = 30
$wsLocal.Cells.Item(8, 8).Value = "https://show.bilibili.com/platform/detail.html?id=82858"
$wsLocal.Cells.Item(8, 9).Value = "//i1.hdslb.com/bfs/openplatform/202403/eeG6Usri1710399065622.jpeg"

# Match the bold/centered/bordered style used by the other column-A index cells
$wsLocal.Range("A7").Copy()
$wsLocal.Range("A8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet "全部类型" (All types) - column F increments
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    3  = 1683
    4  = 459
    5  = 1035
    7  = 969
    9  = 1190
    10 = 948
    12 = 36
    17 = 1051
    19 = 2513
    20 = 576
    22 = 1697
    24 = 648
    28 = 1528
    31 = 772
    32 = 657
    33 = 515
    36 = 45
    39 = 263
    41 = 1165
    42 = 333
    43 = 2462
    47 = 470
    50 = 4036
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Cells.Item($row, 6).Value = $allUpdates[$row]
}
